$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.22"
$ws.Range("E2").Value = "'1.89%"
$ws.Range("D3").Value = "'31.61"
$ws.Range("E3").Value = "'-0.25%"
$ws.Range("D4").Value = "'5.165"
$ws.Range("E4").Value = "'2.98%"
$ws.Range("D5").Value = "'0.07492"
$ws.Range("E5").Value = "'-0.34%"
$ws.Range("D6").Value = "'2.400"
$ws.Range("E6").Value = "'38.93%"
$ws.Range("D7").Value = "'8.015"
$ws.Range("E7").Value = "'2.93%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'3.867"
$ws.Range("E8").Value = "'1.98%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9155"
$ws.Range("E9").Value = "'-1.25%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1733"
$ws.Range("E10").Value = "'1.82%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.07701"
$ws.Range("E11").Value = "'4.86%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08168"
$ws.Range("E12").Value = "'2.97%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03038"
$ws.Range("E13").Value = "'0.00%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09939"
$ws.Range("E14").Value = "'0.39%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001506"
$ws.Range("E15").Value = "'1.16%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006142"
$ws.Range("E16").Value = "'-3.04%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.496"
$ws.Range("E17").Value = "'1.27%"
$ws.Range("D18").Value = "'2.236"
$ws.Range("E18").Value = "'0.89%"
$ws.Range("E19").Value = "'-0.88%"
$ws.Range("E20").Value = "'-0.35%"
$ws.Range("D21").Value = "'4.659"
$ws.Range("E21").Value = "'2.01%"
$ws.Range("D22").Value = "'0.04611"
$ws.Range("E22").Value = "'-0.79%"
$ws.Range("D23").Value = "'0.1564"
$ws.Range("E23").Value = "'0.89%"
$ws.Range("E24").Value = "'3.83%"
$ws.Range("D25").Value = "'0.004533"
$ws.Range("E25").Value = "'2.48%"
$ws.Range("E26").Value = "'-7.28%"
$ws.Range("D27").Value = "'0.0002739"
$ws.Range("E27").Value = "'48.59%"
$ws.Range("D39").Value = "'0.01751"
$ws.Range("E39").Value = "'4.81%"
$ws.Range("D40").Value = "'0.04544"
$ws.Range("E40").Value = "'-0.13%"
$ws.Range("D41").Value = "'0.007421"
$ws.Range("E41").Value = "'5.57%"
$ws.Range("D42").Value = "'0.1363"
$ws.Range("E43").Value = "'5.17%"
$ws.Range("D44").Value = "'0.01090"
$ws.Range("E44").Value = "'-14.84%"
$ws.Range("D45").Value = "'0.00006326"
$ws.Range("E45").Value = "'4.05%"
$ws.Range("E46").Value = "'-57.22%"
